# Insert a new data row at row 3 (pushes existing rows 3-20 down to 4-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new record
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44921
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103001
$ws.Range("J3").Value = "Cereza"
$ws.Range("K3").Value = "Bing"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 320
$ws.Range("N3").Value = 7500
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 7781
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 778
$ws.Range("T3").Value = 10
